$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data added a new weekly price observation. It sits chronologically
# before the current row 23, so a new row is inserted at row 23 and the
# existing rows 23-28 shift down to 24-29 (dimension grows from R28 to R29).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new observation's data.
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44524
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112022
$ws.Range("G23").Value = "Arveja Verde"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12500
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 500
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
